$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Qty executed upto date (column C) updates ---
$ws.Range("C8").Value = 24
$ws.Range("C9").Value = 94
$ws.Range("C10").Value = 54
$ws.Range("C11").Value = 14
$ws.Range("C13").Value = 20
$ws.Range("C14").Value = 67
$ws.Range("C15").Value = 33
$ws.Range("C16").Value = 50
$ws.Range("C17").Value = 24

# --- Upto date Amount (column G) recalculated values, stored as text like the rest of the sheet ---
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "24064.00"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "25488.00"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "9268.00"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "2720.00"

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "1541.00"

# --- Grand totals (rows 19 and 21), columns G and H ---
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "63081.00"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "63081.00"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "63081.00"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "63081.00"
